# The commit swaps the presentation's theme color scheme back to the
# stock "Office Theme" colors (away from the green "Integral" scheme
# that the deck's single slide master/theme currently uses).
#
# Apply the 12-slot Office-Theme color scheme (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) through the live ThemeColorScheme so the
# theme part backing the slide master (and therefore every slide) is
# rewritten in place, matching the target OOXML.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$tcs = $slide.ThemeColorScheme

# MsoThemeColorSchemeIndex order: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1..accent6, 11 hlink, 12 folHlink.
$officeTheme = @(
    @(0,0,0),        # 1  dk1      000000
    @(255,255,255),  # 2  lt1      FFFFFF
    @(68,84,106),    # 3  dk2      44546A
    @(231,230,230),  # 4  lt2      E7E6E6
    @(91,155,213),   # 5  accent1  5B9BD5
    @(237,125,49),   # 6  accent2  ED7D31
    @(165,165,165),  # 7  accent3  A5A5A5
    @(255,192,0),    # 8  accent4  FFC000
    @(68,114,196),   # 9  accent5  4472C4
    @(112,173,71),   # 10 accent6  70AD47
    @(5,99,193),     # 11 hlink    0563C1
    @(149,79,114)    # 12 folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $rgb = $officeTheme[$i - 1]
    $r = $rgb[0]
    $g = $rgb[1]
    $b = $rgb[2]
    $packed = $r + ($g * 256) + ($b * 65536)
    $tcs.Colors($i).RGB = $packed
}
